$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 5.212075734138489
$arr[0,1] = -2.939898788928986
$arr[0,2] = 1.854160755872726
$arr[0,3] = -0.6933320760726929
$arr[0,4] = 1.070083141326904
$arr[0,5] = -0.4738787114620209
$ws.Range("C2:H2").Value = $arr

$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.245875406265258
$arr[0,1] = -4.094250345230103
$arr[0,2] = 2.58136396408081
$arr[0,3] = -0.694248378276825
$arr[0,4] = 0.7021896243095398
$arr[0,5] = -0.1614211350679397
$ws.Range("C3:H3").Value = $arr

$ws.Range("A4").Value = 200
$ws.Range("B4").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.096780717372894
$arr[0,1] = -3.371150910854339
$arr[0,2] = 2.613386332988739
$arr[0,3] = 0.1959350258111953
$arr[0,4] = 0.2964223623275757
$arr[0,5] = -0.2142609804868698
$ws.Range("C4:H4").Value = $arr

$ws.Range("A5").Value = 300
$ws.Range("B5").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.924881196022033
$arr[0,1] = -2.50173692703247
$arr[0,2] = 2.383840799331665
$arr[0,3] = 0.4120286107063293
$arr[0,4] = -0.2721404731273651
$arr[0,5] = -0.3208569586277008
$ws.Range("C5:H5").Value = $arr

$ws.Range("A6").Value = 400
$ws.Range("B6").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.458236134052277
$arr[0,1] = -2.494953083992004
$arr[0,2] = 2.681476718187333
$arr[0,3] = 0.1794416606426239
$arr[0,4] = 0.0216857157647609
$arr[0,5] = -0.3880521357059479
$ws.Range("C6:H6").Value = $arr

$ws.Range("A7").Value = 500
$ws.Range("B7").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.869051647186279
$arr[0,1] = -3.347217082977296
$arr[0,2] = 2.930787801742554
$arr[0,3] = -0.1267545372247696
$arr[0,4] = -0.0375682115554809
$arr[0,5] = -0.1798998117446899
$ws.Range("C7:H7").Value = $arr

$ws.Range("A8").Value = 600
$ws.Range("B8").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.807376968860626
$arr[0,1] = -3.875039219856262
$arr[0,2] = 3.420289939641953
$arr[0,3] = -0.0123700210824608
$arr[0,4] = 0.0419969856739044
$arr[0,5] = 0.271224170923233
$ws.Range("C8:H8").Value = $arr

$ws.Range("A9").Value = 700
$ws.Range("B9").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.295876741409302
$arr[0,1] = -4.034408569335938
$arr[0,2] = 3.253981232643127
$arr[0,3] = -0.0514653958380222
$arr[0,4] = -0.052381694316864
$arr[0,5] = 0.3119994103908539
$ws.Range("C9:H9").Value = $arr

$ws.Range("A10").Value = 800
$ws.Range("B10").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.279258108139039
$arr[0,1] = -3.775099605321884
$arr[0,2] = 3.11082683801651
$arr[0,3] = -0.1554652005434036
$arr[0,4] = -0.0441350154578685
$arr[0,5] = -0.0074830991216003
$ws.Range("C10:H10").Value = $arr

$ws.Range("A11").Value = 900
$ws.Range("B11").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.566667938232422
$arr[0,1] = -3.378203916549682
$arr[0,2] = 3.007539582252503
$arr[0,3] = -0.2000583708286285
$arr[0,4] = -0.1212567538022995
$arr[0,5] = -0.0207694191485643
$ws.Range("C11:H11").Value = $arr

$ws.Range("A12").Value = 1000
$ws.Range("B12").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.106618106365205
$arr[0,1] = -3.249815458059311
$arr[0,2] = 3.031012719869614
$arr[0,3] = -0.1815796941518783
$arr[0,4] = -0.0572686158120632
$arr[0,5] = 0.08643743395805351
$ws.Range("C12:H12").Value = $arr

$ws.Range("A13").Value = 1100
$ws.Range("B13").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.987140679359436
$arr[0,1] = -3.142817544937134
$arr[0,2] = 3.183629143238068
$arr[0,3] = -0.0739146918058395
$arr[0,4] = -0.1140790879726409
$arr[0,5] = 0.1067487001419067
$ws.Range("C13:H13").Value = $arr

$ws.Range("A14").Value = 1200
$ws.Range("B14").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.434188187122345
$arr[0,1] = -3.181812554597855
$arr[0,2] = 3.162444919347763
$arr[0,3] = -0.0395535230636596
$arr[0,4] = -0.0899499058723449
$arr[0,5] = -0.0404698215425014
$ws.Range("C14:H14").Value = $arr

$ws.Range("A15").Value = 1300
$ws.Range("B15").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.282221984863281
$arr[0,1] = -3.265003252029419
$arr[0,2] = 3.094355344772339
$arr[0,3] = -0.0148134818300604
$arr[0,4] = 0.1036943718791008
$arr[0,5] = -0.1157589629292488
$ws.Range("C15:H15").Value = $arr

$ws.Range("A16").Value = 1400
$ws.Range("B16").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.110153055191039
$arr[0,1] = -3.195758980512619
$arr[0,2] = 3.138975620269776
$arr[0,3] = 0.5971207618713379
$arr[0,4] = 1.289536476135254
$arr[0,5] = -0.3637702465057373
$ws.Range("C16:H16").Value = $arr

$ws.Range("A17").Value = 1500
$ws.Range("B17").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.555334329605102
$arr[0,1] = -2.938729083538055
$arr[0,2] = 3.47747951745987
$arr[0,3] = 1.519069194793701
$arr[0,4] = -0.4518875777721405
$arr[0,5] = -0.6734789609909058
$ws.Range("C17:H17").Value = $arr

$ws.Range("A18").Value = 1600
$ws.Range("B18").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.6493126988410929
$arr[0,1] = -2.875420850515366
$arr[0,2] = 3.464587104320525
$arr[0,3] = 0.2113593816757202
$arr[0,4] = -0.3769038617610931
$arr[0,5] = 0.4825835525989532
$ws.Range("C18:H18").Value = $arr

$ws.Range("A19").Value = 1700
$ws.Range("B19").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = -0.6712930202484151
$arr[0,1] = -3.392556905746461
$arr[0,2] = 2.365111112594603
$arr[0,3] = 0.2393064647912979
$arr[0,4] = -0.8791878223419189
$arr[0,5] = -0.1872301995754242
$ws.Range("C19:H19").Value = $arr

$ws.Range("A20").Value = 1800
$ws.Range("B20").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.5537151455879301
$arr[0,1] = -4.361428594589236
$arr[0,2] = 3.347476267814645
$arr[0,3] = -0.1922698318958282
$arr[0,4] = -0.9285151958465576
$arr[0,5] = 0.8594874143600464
$ws.Range("C20:H20").Value = $arr

$ws.Range("A21").Value = 1900
$ws.Range("B21").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.452674245834349
$arr[0,1] = -4.409869003295896
$arr[0,2] = 5.084140586853025
$arr[0,3] = -3.570354700088501
$arr[0,4] = -0.7802276611328125
$arr[0,5] = -4.989242076873779
$ws.Range("C21:H21").Value = $arr

$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 5.681596696376824
$arr[0,1] = -3.525700151920317
$arr[0,2] = 2.85166837722062
$arr[0,3] = -1.221577763557434
$arr[0,4] = 2.375196695327759
$arr[0,5] = -2.503631114959717
$ws.Range("C22:H22").Value = $arr

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.737989616393984
$arr[0,1] = -2.796510410308835
$arr[0,2] = 1.393881118297584
$arr[0,3] = 2.165364503860474
$arr[0,4] = 0.5566509366035461
$arr[0,5] = -0.4453207552433014
$ws.Range("C23:H23").Value = $arr

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = -5.58917605876923
$arr[0,1] = -7.849099040031435
$arr[0,2] = 6.400659620761871
$arr[0,3] = 0.2727513313293457
$arr[0,4] = 0.5925393104553223
$arr[0,5] = 0.4948008358478546
$ws.Range("C24:H24").Value = $arr

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.166972637176496
$arr[0,1] = -11.98566874265667
$arr[0,2] = 9.183138275146462
$arr[0,3] = -0.2121229618787765
$arr[0,4] = 1.80510675907135
$arr[0,5] = 1.96942949295044
$ws.Range("C25:H25").Value = $arr

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = -0.1177038192748863
$arr[0,1] = -2.770210593938835
$arr[0,2] = 3.860614097118379
$arr[0,3] = 0.2755002379417419
$arr[0,4] = 1.588096976280212
$arr[0,5] = 2.037540912628174
$ws.Range("C26:H26").Value = $arr

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.127950906753536
$arr[0,1] = -4.12096252441406
$arr[0,2] = 3.36216964721679
$arr[0,3] = 0.8869763612747192
$arr[0,4] = 0.8231409192085266
$arr[0,5] = 1.362993121147156
$ws.Range("C27:H27").Value = $arr

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 2.318384975194932
$arr[0,1] = -3.449181020259855
$arr[0,2] = 1.645497059822083
$arr[0,3] = 0.2379320114850998
$arr[0,4] = -0.7533495426177979
$arr[0,5] = 0.1786780804395675
$ws.Range("C28:H28").Value = $arr

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.076034724712372
$arr[0,1] = -2.935223340988159
$arr[0,2] = 1.488467574119568
$arr[0,3] = 0.6151412725448608
$arr[0,4] = 1.230893492698669
$arr[0,5] = -0.3686571717262268
$ws.Range("C29:H29").Value = $arr

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 4.307219874858854
$arr[0,1] = -2.422872281074523
$arr[0,2] = 1.113696080446245
$arr[0,3] = 0.1403462886810302
$arr[0,4] = 0.7915286421775818
$arr[0,5] = 0.00137444678694
$ws.Range("C30:H30").Value = $arr

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "falling"
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 3.889325714111331
$arr[0,1] = -2.331348705291748
$arr[0,2] = 1.212168788909908
$arr[0,3] = -0.3019201457500458
$arr[0,4] = 0.041233405470848
$arr[0,5] = -0.0345138870179653
$ws.Range("C31:H31").Value = $arr
